# "Fruta / hortaliza, semanal"
# A new weekly price observation is inserted into the daily log as row 28,
# pushing the existing rows 28-56 down to 29-57 (same data, just shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28, shifting rows 28:56 down to 29:57.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly observation.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44586
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 42000
$ws.Range("L28").Value = 42000
$ws.Range("M28").Value = 42000
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1680
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
